# Fruta / hortaliza, semanal
# Insert two new weekly price rows for "Femacal de La Calera" / Frutilla
# just above the existing row that will become row 154 (old row 152).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 152:153 - this shifts old rows 152-160 down to 154-162,
# carrying their formatting/content with them automatically.
$ws.Rows("152:153").Insert()

# New row 152 - "Especial" quality, week of 2021-11-16
$ws.Cells.Item(152, 1).Value = 3
$ws.Cells.Item(152, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(152, 3).Value = "Coquimbo"
$ws.Cells.Item(152, 4).Value = 44516
$ws.Cells.Item(152, 5).Value = 5
$ws.Cells.Item(152, 6).Value = "Fruta"
$ws.Cells.Item(152, 7).Value = 100101
$ws.Cells.Item(152, 8).Value = "Berries"
$ws.Cells.Item(152, 9).Value = 100112025
$ws.Cells.Item(152, 10).Value = "Frutilla"
$ws.Cells.Item(152, 11).Value = "Sin especificar"
$ws.Cells.Item(152, 12).Value = "Especial"
$ws.Cells.Item(152, 13).Value = 65
$ws.Cells.Item(152, 14).Value = 6000
$ws.Cells.Item(152, 15).Value = 6000
$ws.Cells.Item(152, 16).Value = 6000
$ws.Cells.Item(152, 17).Value = "`$/bandeja 7 kilos"
$ws.Cells.Item(152, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(152, 19).Value = 857
$ws.Cells.Item(152, 20).Value = 7

# New row 153 - "Segunda" quality, week of 2021-11-16
$ws.Cells.Item(153, 1).Value = 3
$ws.Cells.Item(153, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(153, 3).Value = "Coquimbo"
$ws.Cells.Item(153, 4).Value = 44516
$ws.Cells.Item(153, 5).Value = 5
$ws.Cells.Item(153, 6).Value = "Fruta"
$ws.Cells.Item(153, 7).Value = 100101
$ws.Cells.Item(153, 8).Value = "Berries"
$ws.Cells.Item(153, 9).Value = 100112025
$ws.Cells.Item(153, 10).Value = "Frutilla"
$ws.Cells.Item(153, 11).Value = "Sin especificar"
$ws.Cells.Item(153, 12).Value = "Segunda"
$ws.Cells.Item(153, 13).Value = 60
$ws.Cells.Item(153, 14).Value = 4000
$ws.Cells.Item(153, 15).Value = 4000
$ws.Cells.Item(153, 16).Value = 4000
$ws.Cells.Item(153, 17).Value = "`$/bandeja 7 kilos"
$ws.Cells.Item(153, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(153, 19).Value = 571
$ws.Cells.Item(153, 20).Value = 7
